$wb = $excel.ActiveWorkbook

# OFF sheet (sheet1) - row 2 values
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 657
$wsOff.Range("C2").Value = 458
$wsOff.Range("D2").Value = 154
$wsOff.Range("E2").Value = 75

# DEF sheet (sheet2) - row 2 values
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 607
$wsDef.Range("C2").Value = 423
$wsDef.Range("D2").Value = 138
$wsDef.Range("E2").Value = 50
